$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 21.18599966666667
$ws.Range("H2").Value = 63.557999
$ws.Range("I2").Value = 0.08765141600314529
$ws.Range("J2").Value = 0.08765141600314529
$ws.Range("M2").Value = 1.758668
$ws.Range("N2").Value = 5.276004
$ws.Range("O2").Value = 0.02465283256602696
$ws.Range("P2").Value = 0.02465283256602696
$ws.Range("Q2").Value = 37.25913966177734
$ws.Range("R2").Value = 335.3322569559961
$ws.Range("S2").Value = 0.002160855682900717
$ws.Range("T2").Value = 0.002160855682900717
$ws.Range("G3").Value = 21.18599966666667
$ws.Range("H3").Value = 63.557999
$ws.Range("I3").Value = 0.08765141600314529
$ws.Range("J3").Value = 0.08765141600314529
$ws.Range("O3").Value = 0.0796780206066965
$ws.Range("P3").Value = 0.0796780206066965
$ws.Range("Q3").Value = 120.4216387633268
$ws.Range("R3").Value = 1083.794748869941
$ws.Range("S3").Value = 0.006983891330504737
$ws.Range("T3").Value = 0.006983891330504737
$ws.Range("G4").Value = 21.18599966666667
$ws.Range("H4").Value = 63.557999
$ws.Range("I4").Value = 0.08765141600314529
$ws.Range("J4").Value = 0.08765141600314529
$ws.Range("M4").Value = 3.568404
$ws.Range("N4").Value = 10.705212
$ws.Range("O4").Value = 0.05002153126112539
$ws.Range("P4").Value = 0.05002153126112539
$ws.Range("Q4").Value = 75.600205954532
$ws.Range("R4").Value = 680.4018535907881
$ws.Range("S4").Value = 0.004384458045683238
$ws.Range("T4").Value = 0.004384458045683238
$ws.Range("G5").Value = 21.18599966666667
$ws.Range("H5").Value = 63.557999
$ws.Range("I5").Value = 0.08765141600314529
$ws.Range("J5").Value = 0.08765141600314529
$ws.Range("M5").Value = 60.32626866666666
$ws.Range("N5").Value = 180.978806
$ws.Range("O5").Value = 0.8456476155661511
$ws.Range("P5").Value = 0.8456476155661511
$ws.Range("Q5").Value = 1278.072307863244
$ws.Range("R5").Value = 11502.65077076919
$ws.Range("S5").Value = 0.07412221094405659
$ws.Range("T5").Value = 0.07412221094405659
$ws.Range("I6").Value = 0.5040014103551328
$ws.Range("J6").Value = 0.5040014103551328
$ws.Range("M6").Value = 1.758668
$ws.Range("N6").Value = 5.276004
$ws.Range("O6").Value = 0.02465283256602696
$ws.Range("P6").Value = 0.02465283256602696
$ws.Range("Q6").Value = 214.2425050780787
$ws.Range("R6").Value = 1928.182545702708
$ws.Range("S6").Value = 0.01242506238252654
$ws.Range("T6").Value = 0.01242506238252654
$ws.Range("I7").Value = 0.5040014103551328
$ws.Range("J7").Value = 0.5040014103551328
$ws.Range("O7").Value = 0.0796780206066965
$ws.Range("P7").Value = 0.0796780206066965
$ws.Range("S7").Value = 0.04015783476008037
$ws.Range("T7").Value = 0.04015783476008037
$ws.Range("I8").Value = 0.5040014103551328
$ws.Range("J8").Value = 0.5040014103551328
$ws.Range("M8").Value = 3.568404
$ws.Range("N8").Value = 10.705212
$ws.Range("O8").Value = 0.05002153126112539
$ws.Range("P8").Value = 0.05002153126112539
$ws.Range("Q8").Value = 434.706159485836
$ws.Range("R8").Value = 3912.355435372524
$ws.Range("S8").Value = 0.02521092230373056
$ws.Range("T8").Value = 0.02521092230373056
$ws.Range("I9").Value = 0.5040014103551328
$ws.Range("J9").Value = 0.5040014103551328
$ws.Range("M9").Value = 60.32626866666666
$ws.Range("N9").Value = 180.978806
$ws.Range("O9").Value = 0.8456476155661511
$ws.Range("P9").Value = 0.8456476155661511
$ws.Range("Q9").Value = 7348.999880113741
$ws.Range("R9").Value = 66140.99892102367
$ws.Range("S9").Value = 0.4262075909087953
$ws.Range("T9").Value = 0.4262075909087953
$ws.Range("G10").Value = 37.20718233333333
$ws.Range("H10").Value = 111.621547
$ws.Range("I10").Value = 0.1539347809079331
$ws.Range("J10").Value = 0.1539347809079331
$ws.Range("M10").Value = 1.758668
$ws.Range("N10").Value = 5.276004
$ws.Range("O10").Value = 0.02465283256602696
$ws.Range("P10").Value = 0.02465283256602696
$ws.Range("Q10").Value = 65.43508093979867
$ws.Range("R10").Value = 588.915728458188
$ws.Range("S10").Value = 0.003794928379811319
$ws.Range("T10").Value = 0.00379492837981132
$ws.Range("G11").Value = 37.20718233333333
$ws.Range("H11").Value = 111.621547
$ws.Range("I11").Value = 0.1539347809079331
$ws.Range("J11").Value = 0.1539347809079331
$ws.Range("O11").Value = 0.0796780206066965
$ws.Range("P11").Value = 0.0796780206066965
$ws.Range("Q11").Value = 211.4863561239192
$ws.Range("R11").Value = 1903.377205115273
$ws.Range("S11").Value = 0.0122652186452696
$ws.Range("T11").Value = 0.01226521864526961
$ws.Range("G12").Value = 37.20718233333333
$ws.Range("H12").Value = 111.621547
$ws.Range("I12").Value = 0.1539347809079331
$ws.Range("J12").Value = 0.1539347809079331
$ws.Range("M12").Value = 3.568404
$ws.Range("N12").Value = 10.705212
$ws.Range("O12").Value = 0.05002153126112539
$ws.Range("P12").Value = 0.05002153126112539
$ws.Range("Q12").Value = 132.770258266996
$ws.Range("R12").Value = 1194.932324402964
$ws.Range("S12").Value = 0.007700053455360663
$ws.Range("T12").Value = 0.007700053455360666
$ws.Range("G13").Value = 37.20718233333333
$ws.Range("H13").Value = 111.621547
$ws.Range("I13").Value = 0.1539347809079331
$ws.Range("J13").Value = 0.1539347809079331
$ws.Range("M13").Value = 60.32626866666666
$ws.Range("N13").Value = 180.978806
$ws.Range("O13").Value = 0.8456476155661511
$ws.Range("P13").Value = 0.8456476155661511
$ws.Range("Q13").Value = 2244.57047777032
$ws.Range("R13").Value = 20201.13429993288
$ws.Range("S13").Value = 0.1301745804274915
$ws.Range("T13").Value = 0.1301745804274915
$ws.Range("G14").Value = 61.49336899999999
$ws.Range("H14").Value = 184.480107
$ws.Range("I14").Value = 0.2544123927337887
$ws.Range("J14").Value = 0.2544123927337887
$ws.Range("M14").Value = 1.758668
$ws.Range("N14").Value = 5.276004
$ws.Range("O14").Value = 0.02465283256602696
$ws.Range("P14").Value = 0.02465283256602696
$ws.Range("Q14").Value = 108.146420272492
$ws.Range("R14").Value = 973.317782452428
$ws.Range("S14").Value = 0.006271986120788388
$ws.Range("T14").Value = 0.006271986120788389
$ws.Range("G15").Value = 61.49336899999999
$ws.Range("H15").Value = 184.480107
$ws.Range("I15").Value = 0.2544123927337887
$ws.Range("J15").Value = 0.2544123927337887
$ws.Range("O15").Value = 0.0796780206066965
$ws.Range("P15").Value = 0.0796780206066965
$ws.Range("Q15").Value = 349.5295187655903
$ws.Range("R15").Value = 3145.765668890313
$ws.Range("S15").Value = 0.02027107587084178
$ws.Range("T15").Value = 0.02027107587084178
$ws.Range("G16").Value = 61.49336899999999
$ws.Range("H16").Value = 184.480107
$ws.Range("I16").Value = 0.2544123927337887
$ws.Range("J16").Value = 0.2544123927337887
$ws.Range("M16").Value = 3.568404
$ws.Range("N16").Value = 10.705212
$ws.Range("O16").Value = 0.05002153126112539
$ws.Range("P16").Value = 0.05002153126112539
$ws.Range("Q16").Value = 219.433183913076
$ws.Range("R16").Value = 1974.898655217684
$ws.Range("S16").Value = 0.01272609745635092
$ws.Range("T16").Value = 0.01272609745635092
$ws.Range("G17").Value = 61.49336899999999
$ws.Range("H17").Value = 184.480107
$ws.Range("I17").Value = 0.2544123927337887
$ws.Range("J17").Value = 0.2544123927337887
$ws.Range("M17").Value = 60.32626866666666
$ws.Range("N17").Value = 180.978806
$ws.Range("O17").Value = 0.8456476155661511
$ws.Range("P17").Value = 0.8456476155661511
$ws.Range("Q17").Value = 3709.665499512471
$ws.Range("R17").Value = 33386.98949561224
$ws.Range("S17").Value = 0.2151432332858076
$ws.Range("T17").Value = 0.2151432332858076
